# Updates cryptos list prices/volumes (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Preserve the original "Text" cell type (as in the source files inline
    # strings) instead of letting Excel auto-coerce numeric-looking strings
    # (e.g. "0.9989") into actual numbers. Restores the prior cell style
    # afterwards so no visible formatting changes.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "25.659.83"
$ws.Range("E2").Value = "  +5.03%  "
Set-TextValue $ws.Range("D3") "1.706.48"
$ws.Range("E3").Value = "  +3.26%  "
Set-TextValue $ws.Range("D4") "0.9989"
$ws.Range("E4").Value = "  -0.53%  "
Set-TextValue $ws.Range("D5") "330.69"
$ws.Range("E5").Value = "  +5.88%  "
Set-TextValue $ws.Range("D6") "0.9965"
$ws.Range("E6").Value = "  -0.53%  "
Set-TextValue $ws.Range("D7") "0.3670"
$ws.Range("E7").Value = "  +0.90%  "
Set-TextValue $ws.Range("D8") "48.42"
$ws.Range("E8").Value = "  +3.08%  "
Set-TextValue $ws.Range("D9") "0.3286"
$ws.Range("E9").Value = "  +0.91%  "
Set-TextValue $ws.Range("D10") "1.164"
$ws.Range("E10").Value = "  +3.70%  "
Set-TextValue $ws.Range("D11") "0.07312"
$ws.Range("E11").Value = "  +3.83%  "
Set-TextValue $ws.Range("D12") "0.9970"
$ws.Range("E12").Value = "  -0.41%  "
Set-TextValue $ws.Range("D13") "6.183"
$ws.Range("E13").Value = "  +4.16%  "
Set-TextValue $ws.Range("D14") "19.91"
$ws.Range("E14").Value = "  +2.50%  "
Set-TextValue $ws.Range("D15") "1.701.38"
$ws.Range("E15").Value = "  +2.96%  "
Set-TextValue $ws.Range("D16") "6.797"
$ws.Range("E16").Value = "  +3.13%  "
Set-TextValue $ws.Range("D17") "0.00001068"
$ws.Range("E17").Value = "  +2.21%  "
Set-TextValue $ws.Range("D18") "0.06586"
$ws.Range("E18").Value = "  -0.35%  "
Set-TextValue $ws.Range("D19") "80.76"
$ws.Range("E19").Value = "  +3.44%  "
Set-TextValue $ws.Range("D20") "0.9957"
$ws.Range("E20").Value = "  -0.52%  "
Set-TextValue $ws.Range("D21") "6.035"
$ws.Range("E21").Value = "  +1.82%  "
Set-TextValue $ws.Range("D22") "16.10"
$ws.Range("E22").Value = "  +2.85%  "
Set-TextValue $ws.Range("D23") "12.99"
$ws.Range("E23").Value = "  +4.65%  "
Set-TextValue $ws.Range("D24") "25.634.37"
$ws.Range("E24").Value = "  +5.08%  "
Set-TextValue $ws.Range("D25") "2.448"
$ws.Range("E25").Value = "  -0.70%  "
Set-TextValue $ws.Range("D26") "2.474"
$ws.Range("E26").Value = "  +5.83%  "
Set-TextValue $ws.Range("D27") "149.53"
$ws.Range("E27").Value = "  +0.87%  "
Set-TextValue $ws.Range("D28") "19.07"
$ws.Range("E28").Value = "  +2.55%  "
Set-TextValue $ws.Range("D29") "1.271"
$ws.Range("E29").Value = "  +7.60%  "
Set-TextValue $ws.Range("D30") "1.890.52"
$ws.Range("E30").Value = "  +2.79%  "
Set-TextValue $ws.Range("D31") "128.00"
$ws.Range("E31").Value = "  +3.07%  "
Set-TextValue $ws.Range("D32") "4.106"
$ws.Range("E32").Value = "  +0.40%  "
Set-TextValue $ws.Range("D33") "5.953"
$ws.Range("E33").Value = "  +4.92%  "
Set-TextValue $ws.Range("D34") "0.08495"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("E36").Value = "  +3.02%  "
Set-TextValue $ws.Range("D37") "5.287"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("E38").Value = "  +5.88%  "
Set-TextValue $ws.Range("D39") "0.06196"
$ws.Range("E39").Value = "  +2.65%  "
Set-TextValue $ws.Range("D40") "0.2117"
$ws.Range("E40").Value = "  +2.40%  "
Set-TextValue $ws.Range("D41") "0.02258"
$ws.Range("E41").Value = "  +1.94%  "
Set-TextValue $ws.Range("D42") "8.450"
$ws.Range("E42").Value = "  +3.47%  "
Set-TextValue $ws.Range("D43") "0.6079"
$ws.Range("E43").Value = "  +2.97%  "
Set-TextValue $ws.Range("D44") "0.9961"
$ws.Range("E44").Value = "  -0.46%  "
Set-TextValue $ws.Range("D45") "13.99"
$ws.Range("E45").Value = "  +10.88%  "
Set-TextValue $ws.Range("D46") "3.840"
$ws.Range("E46").Value = "  +1.49%  "
Set-TextValue $ws.Range("D47") "0.5835"
$ws.Range("E47").Value = "  +3.73%  "
Set-TextValue $ws.Range("D48") "125.61"
$ws.Range("E48").Value = "  +2.65%  "
Set-TextValue $ws.Range("D49") "1.996"
$ws.Range("E49").Value = "  +2.78%  "
Set-TextValue $ws.Range("D50") "0.07222"
$ws.Range("E50").Value = "  +4.74%  "
Set-TextValue $ws.Range("D51") "1.207"
$ws.Range("E51").Value = "  +3.99%  "
